# Loan RBI, Variable Instalments
# Insert a new (empty) column before column N on the "Repayment Schedule"
# sheet - this shifts the old N/O/P columns (and their header labels /
# values) one column to the right, and adds a blank N column throughout.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")
$ws.Columns("N").Insert()

# Make "Repayment Schedule" the active sheet/tab, with L18 selected,
# instead of "NewLoanInput".
$ws.Activate()
$ws.Range("L18").Select() | Out-Null
